$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Work from the bottom of the document upward so that paragraph-count
# changing edits (paragraph insertions) never invalidate the indices of
# edits still to be performed further up the document.
# ---------------------------------------------------------------------------

# --- 6) empty paragraph (0-based #48 / 1-based #49, right after
#        "Strømgaten herre:1 pris:10 åpen:nå") gets a new sentence of text.
$p = $d.Paragraphs(49)
$r = $p.Range
$r.End = $r.End - 1
$r.InsertAfter("Det er også mulig å søke fritekst på engelsk på plasseringer.")
$r.Font.Size = 12

# --- 5) remove the stray "_GoBack" bookmark that used to sit at the end of
#        the "... I så fall må man skrive plassering som et kriterie som
#        vist over." paragraph (0-based #46 / 1-based #47).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 4) "Åpen data ting.html ??" section becomes "Utveksling.html" and its
#        body paragraph is rewritten; a brand-new paragraph with the
#        "Land -> dropdown ..." note is added right after the rewritten
#        body paragraph.
$d.Content.Find.Execute("Åpne data ting og tang", $true, $false, $false, `
    $false, $false, $true, 1, $false, `
    "På dette dokumentet har vi importert et json dokument som inneholder informasjon om alle utvekslingsavtalene til universitetet i Bergen. Her er det mulig å søke etter land for å se om UiB har utvekslingsavtale og i så fall hvor mange avtaler som finnes med ulike universiteter.", `
    2)

$p = $d.Paragraphs(39)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertParagraphAfter()
$newp = $d.Paragraphs(40)
$nr = $newp.Range
$nr.End = $nr.End - 1
$nr.InsertAfter("Land -> dropdown til fakultet/institutt?")
$nr.Font.Size = 12

# --- 4b) rename the heading itself (also drops the stale
#         lastRenderedPageBreak cache marker on that run).
$d.Content.Find.Execute("Åpen data ting.html ??", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Utveksling.html", 2)

# --- 3) append a new sentence to the "... også på kartet." paragraph
#        (0-based #32 / 1-based #33).
$p = $d.Paragraphs(33)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(" Lekeplassene er nummerert så man kan se på kartet hvor de er posisjonert.")
$r.Font.Size = 12

# --- 2) "Favtoalett.html" heading: no text change, but rewriting via
#        Find/Replace drops the stale lastRenderedPageBreak cache marker.
$d.Content.Find.Execute("Favtoalett.html", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Favtoalett.html", 2)

# --- 1) the empty paragraph right after "... via funksjonen find_common."
#        (0-based #23 / 1-based #24) gets three new runs of text plus a
#        relocated "_GoBack" bookmark, and a new empty paragraph is
#        inserted right after it.
$p = $d.Paragraphs(24)
$r = $p.Range
$r.End = $r.End - 1

$r.InsertAfter("På plassering er ")
$r.Font.Size = 12
$r.Collapse(0)

$r.InsertAfter("det ")
$r.Font.Size = 12
$r.Collapse(0)
$bookmarkPos = $r.Start

$r.InsertAfter("både adresse og plassering søk.")
$r.Font.Size = 12

$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$p.Range.InsertParagraphAfter()

Write-Output "done"
